$wb = $excel.ActiveWorkbook

# Excel constants
$xlNone         = -4142   # xlLineStyleNone
$xlThin         = 1       # xlContinuous
$xlEdgeLeft     = 7
$xlEdgeTop      = 8
$xlEdgeBottom   = 9
$xlEdgeRight    = 10
$xlPasteFormats = -4122   # xlPasteFormats

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---- Build the two new border styles once (on quality_comparison!C1 / D1) ----
# borderId 4 equivalent: top + bottom thin only
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
$c1.Borders.Item($xlEdgeRight).LineStyle = $xlNone
$c1.Borders.Item($xlEdgeTop).LineStyle = $xlThin
$c1.Borders.Item($xlEdgeBottom).LineStyle = $xlThin

# borderId 5 equivalent: top + bottom + right thin
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
$d1.Borders.Item($xlEdgeRight).LineStyle = $xlThin
$d1.Borders.Item($xlEdgeTop).LineStyle = $xlThin
$d1.Borders.Item($xlEdgeBottom).LineStyle = $xlThin

# Reuse those exact formats everywhere else they're needed (copy/paste-format
# avoids re-deriving the same border combination from scratch, which keeps
# every target cell pointing at the same shared style record).
$c1.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$d1.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---- Text updates: "fedcore" -> "approach" ----
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---- Numeric "-0" -> "0" fixups ----
$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("D12").Value = 0

# ---- Remove stray empty inline-string cell G5 on computational_comparison ----
$ws2.Range("G5").ClearContents()
